# Auto-generated edit script applying scheduled market-data refresh
# to the per-sheet Leve profit tables (columns H-N) as captured by the commit diff.
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 367221.44
$ws.Range("J17").Value = 367221.44
$ws.Range("L17").Value = 1101664.32
$ws.Range("N17").Value = -1102000.32
$ws.Range("H98").Value = 1904.5454
$ws.Range("I98").Value = 1920
$ws.Range("J98").Value = 1750
$ws.Range("K98").Value = 1920
$ws.Range("L98").Value = 1750
$ws.Range("M98").Value = -422
$ws.Range("N98").Value = -4746
$ws.Range("H113").Value = 2396.487
$ws.Range("I113").Value = 2616.476
$ws.Range("J113").Value = 2139.8333
$ws.Range("K113").Value = 2616.476
$ws.Range("L113").Value = 2139.8333
$ws.Range("M113").Value = 637.5239999999999
$ws.Range("N113").Value = -8647.8333
$ws.Range("H122").Value = 1904.5454
$ws.Range("I122").Value = 1920
$ws.Range("J122").Value = 1750
$ws.Range("K122").Value = 5760
$ws.Range("L122").Value = 5250
$ws.Range("M122").Value = -3310
$ws.Range("N122").Value = -10150
$ws.Range("H129").Value = 1001925.2
$ws.Range("J129").Value = 1482724.2
$ws.Range("L129").Value = 4448172.6
$ws.Range("N129").Value = -4458172.6
$ws.Range("H132").Value = 1744.6182
$ws.Range("I132").Value = 1667.1063
$ws.Range("J132").Value = 2200
$ws.Range("K132").Value = 5001.3189
$ws.Range("L132").Value = 6600
$ws.Range("M132").Value = -2471.3189
$ws.Range("N132").Value = -11660
$ws.Range("H138").Value = 2374.63
$ws.Range("I138").Value = 844
$ws.Range("J138").Value = 2999.817
$ws.Range("K138").Value = 2532
$ws.Range("L138").Value = 8999.451000000001
$ws.Range("M138").Value = 2608
$ws.Range("N138").Value = -19279.451

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3719396.2
$ws.Range("I32").Value = 4584413
$ws.Range("J32").Value = 4911.5293
$ws.Range("K32").Value = 4584413
$ws.Range("L32").Value = 4911.5293
$ws.Range("M32").Value = -4584126
$ws.Range("N32").Value = -5485.5293
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("M39").ClearContents()
$ws.Range("N39").ClearContents()
$ws.Range("H132").Value = 1240
$ws.Range("I132").Value = 728.67566
$ws.Range("J132").Value = 3604.875
$ws.Range("K132").Value = 2186.02698
$ws.Range("L132").Value = 10814.625
$ws.Range("M132").Value = 343.9730199999999
$ws.Range("N132").Value = -15874.625

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3314.8096
$ws.Range("I105").Value = 4026.25
$ws.Range("J105").Value = 1038.2
$ws.Range("K105").Value = 4026.25
$ws.Range("L105").Value = 1038.2
$ws.Range("M105").Value = -2279.25
$ws.Range("N105").Value = -4532.2

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 4966.28
$ws.Range("I58").Value = 905.2143
$ws.Range("J58").Value = 10134.909
$ws.Range("K58").Value = 905.2143
$ws.Range("L58").Value = 10134.909
$ws.Range("M58").Value = -702.2143
$ws.Range("N58").Value = -10540.909
$ws.Range("H118").Value = 24850
$ws.Range("I118").Value = 10000
$ws.Range("J118").Value = 29800
$ws.Range("K118").Value = 10000
$ws.Range("L118").Value = 29800
$ws.Range("M118").Value = -8343
$ws.Range("N118").Value = -33114
$ws.Range("H122").Value = 1368.4286
$ws.Range("I122").Value = 1225.5
$ws.Range("J122").Value = 1781.3334
$ws.Range("K122").Value = 3676.5
$ws.Range("L122").Value = 5344.0002
$ws.Range("M122").Value = -1226.5
$ws.Range("N122").Value = -10244.0002
$ws.Range("H132").Value = 1571.3334
$ws.Range("I132").Value = 833.8049
$ws.Range("J132").Value = 5891.143
$ws.Range("K132").Value = 2501.4147
$ws.Range("L132").Value = 17673.429
$ws.Range("M132").Value = 28.58530000000019
$ws.Range("N132").Value = -22733.429
$ws.Range("H136").Value = 4966.28
$ws.Range("I136").Value = 905.2143
$ws.Range("J136").Value = 10134.909
$ws.Range("K136").Value = 2715.6429
$ws.Range("L136").Value = 30404.727
$ws.Range("M136").Value = -165.6428999999998
$ws.Range("N136").Value = -35504.727

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 680.7059
$ws.Range("I68").Value = 510.57144
$ws.Range("J68").Value = 799.8
$ws.Range("K68").Value = 1531.71432
$ws.Range("L68").Value = 2399.4
$ws.Range("M68").Value = -720.71432
$ws.Range("N68").Value = -4021.4
$ws.Range("H71").Value = 680.7059
$ws.Range("I71").Value = 510.57144
$ws.Range("J71").Value = 799.8
$ws.Range("K71").Value = 4595.14296
$ws.Range("L71").Value = 7198.2
$ws.Range("M71").Value = -539.1429600000001
$ws.Range("N71").Value = -15310.2
$ws.Range("H80").Value = 5943
$ws.Range("I80").Value = 2061
$ws.Range("J80").Value = 9825
$ws.Range("K80").Value = 6183
$ws.Range("L80").Value = 29475
$ws.Range("M80").Value = -5247
$ws.Range("N80").Value = -31347
$ws.Range("H83").Value = 5943
$ws.Range("I83").Value = 2061
$ws.Range("J83").Value = 9825
$ws.Range("K83").Value = 18549
$ws.Range("L83").Value = 88425
$ws.Range("M83").Value = -13869
$ws.Range("N83").Value = -97785
$ws.Range("H113").Value = 1300
$ws.Range("I113").Value = 1450
$ws.Range("J113").Value = 1150
$ws.Range("K113").Value = 4350
$ws.Range("L113").Value = 3450
$ws.Range("M113").Value = -2180
$ws.Range("N113").Value = -7790
$ws.Range("H131").Value = 9340717
$ws.Range("I131").Value = 91835210
$ws.Range("J131").Value = 1718.6604
$ws.Range("K131").Value = 275505630
$ws.Range("L131").Value = 5155.9812
$ws.Range("M131").Value = -275500590
$ws.Range("N131").Value = -15235.9812
$ws.Range("H132").Value = 1250
$ws.Range("J132").Value = 2000
$ws.Range("L132").Value = 18000
$ws.Range("N132").Value = -23060

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1443.75
$ws.Range("I22").Value = 1466.75
$ws.Range("J22").Value = 1374.75
$ws.Range("K22").Value = 1466.75
$ws.Range("L22").Value = 1374.75
$ws.Range("M22").Value = -1171.75
$ws.Range("N22").Value = -1964.75
$ws.Range("H27").Value = 1443.75
$ws.Range("I27").Value = 1466.75
$ws.Range("J27").Value = 1374.75
$ws.Range("K27").Value = 1466.75
$ws.Range("L27").Value = 1374.75
$ws.Range("M27").Value = -1359.75
$ws.Range("N27").Value = -1588.75

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 11214.143
$ws.Range("J14").Value = 11214.143
$ws.Range("L14").Value = 11214.143
$ws.Range("N14").Value = -11550.143
$ws.Range("H116").Value = 30000
$ws.Range("J116").Value = 30000
$ws.Range("L116").Value = 30000
$ws.Range("N116").Value = -39178

